$wb = $excel.ActiveWorkbook

# Avoid the "are you sure you want to delete" prompt when removing a sheet.
$excel.DisplayAlerts = $false

# Rename sheets (content case / accent updates).
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Remove the obsolete "Desarquivamentos Pendentes" sheet entirely.
[void]$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()

# Deleting a sheet shifts the active tab; restore the original active sheet.
$wb.Worksheets.Item("PAINEIS DARQ").Activate()

$excel.DisplayAlerts = $true
